$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 9938.799999999999
$ws.Range("I32").Value = 4750
$ws.Range("K32").Value = 4750
$ws.Range("M32").Value = -4424

$ws.Range("H51").Value = 4999.7144
$ws.Range("J51").Value = 4999.7144
$ws.Range("L51").Value = 4999.7144
$ws.Range("N51").Value = -5967.7144

$ws.Range("H70").Value = 4317
$ws.Range("J70").Value = 4180.4
$ws.Range("L70").Value = 12541.2
$ws.Range("N70").Value = -13081.2

$ws.Range("H73").Value = 4317
$ws.Range("J73").Value = 4180.4
$ws.Range("L73").Value = 12541.2
$ws.Range("N73").Value = -14413.2

$ws.Range("H80").Value = 2122.963
$ws.Range("I80").Value = 742.4
$ws.Range("J80").Value = 2436.7273
$ws.Range("K80").Value = 2227.2
$ws.Range("L80").Value = 7310.1819
$ws.Range("M80").Value = -1229.2
$ws.Range("N80").Value = -9306.1819

$ws.Range("H83").Value = 2122.963
$ws.Range("I83").Value = 742.4
$ws.Range("J83").Value = 2436.7273
$ws.Range("K83").Value = 6681.599999999999
$ws.Range("L83").Value = 21930.5457
$ws.Range("M83").Value = -1689.599999999999
$ws.Range("N83").Value = -31914.5457

$ws.Range("H88").Value = 8455.571
$ws.Range("I88").Value = 300
$ws.Range("J88").Value = 9814.833000000001
$ws.Range("K88").Value = 300
$ws.Range("L88").Value = 9814.833000000001
$ws.Range("M88").Value = 106
$ws.Range("N88").Value = -10626.833

$ws.Range("H91").Value = 8455.571
$ws.Range("I91").Value = 300
$ws.Range("J91").Value = 9814.833000000001
$ws.Range("K91").Value = 300
$ws.Range("L91").Value = 9814.833000000001
$ws.Range("M91").Value = 1104
$ws.Range("N91").Value = -12622.833

$ws.Range("H135").Value = 2903
$ws.Range("I135").Value = 3028.9375
$ws.Range("J135").Value = 2500
$ws.Range("K135").Value = 27260.4375
$ws.Range("L135").Value = 22500
$ws.Range("M135").Value = -24725.4375
$ws.Range("N135").Value = -27570

$ws.Range("H137").Value = 3820.8096
$ws.Range("I137").Value = 1605.1111
$ws.Range("J137").Value = 5482.5835
$ws.Range("K137").Value = 4815.3333
$ws.Range("L137").Value = 16447.7505
$ws.Range("M137").Value = -2265.3333
$ws.Range("N137").Value = -21547.7505

$ws.Range("H138").Value = 11153.692
$ws.Range("I138").Value = 2999.6
$ws.Range("J138").Value = 16250
$ws.Range("K138").Value = 8998.799999999999
$ws.Range("L138").Value = 48750
$ws.Range("M138").Value = -3858.799999999999
$ws.Range("N138").Value = -59030

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 479.5
$ws.Range("I4").Value = 479.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 479.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -363.5
$ws.Range("N4").ClearContents()

$ws.Range("H5").Value = 606.6
$ws.Range("I5").Value = 606.6
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 606.6
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -494.6
$ws.Range("N5").ClearContents()

$ws.Range("H32").Value = 1504.3334
$ws.Range("I32").Value = 1037.2031
$ws.Range("J32").Value = 4222.1816
$ws.Range("K32").Value = 1037.2031
$ws.Range("L32").Value = 4222.1816
$ws.Range("M32").Value = -750.2030999999999
$ws.Range("N32").Value = -4796.1816

$ws.Range("H36").Value = 18332.666
$ws.Range("I36").Value = 9999
$ws.Range("J36").Value = 22499.5
$ws.Range("K36").Value = 9999
$ws.Range("L36").Value = 22499.5
$ws.Range("N36").Value = -23191.5
$ws.Range("M36").Value = -9653

$ws.Range("H74").Value = 26318018
$ws.Range("I74").Value = 38463212
$ws.Range("J74").Value = 3430.5
$ws.Range("K74").Value = 38463212
$ws.Range("L74").Value = 3430.5
$ws.Range("M74").Value = -38462338
$ws.Range("N74").Value = -5178.5

$ws.Range("H76").Value = 59999.75
$ws.Range("J76").Value = 59999.75
$ws.Range("L76").Value = 59999.75
$ws.Range("N76").Value = -60675.75

$ws.Range("H77").Value = 26318018
$ws.Range("I77").Value = 38463212
$ws.Range("J77").Value = 3430.5
$ws.Range("K77").Value = 192316060
$ws.Range("L77").Value = 17152.5
$ws.Range("M77").Value = -192311692
$ws.Range("N77").Value = -25888.5

$ws.Range("H79").Value = 59999.75
$ws.Range("J79").Value = 59999.75
$ws.Range("L79").Value = 59999.75
$ws.Range("N79").Value = -62339.75

$ws.Range("H97").Value = 1085.4
$ws.Range("I97").Value = 932.9666999999999
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 932.9666999999999
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -436.9666999999999
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 606.6
$ws.Range("I4").Value = 606.6
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 606.6
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -491.6
$ws.Range("N4").ClearContents()

$ws.Range("H94").Value = 853.5217
$ws.Range("I94").Value = 547.44446
$ws.Range("J94").Value = 1955.4
$ws.Range("K94").Value = 547.44446
$ws.Range("L94").Value = 1955.4
$ws.Range("M94").Value = -96.44446000000005
$ws.Range("N94").Value = -2857.4

$ws.Range("H134").Value = 2645.8572
$ws.Range("I134").Value = 2278.15
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 6834.450000000001
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -4299.450000000001
$ws.Range("N134").Value = -35070

$ws.Range("H138").Value = 149993.25
$ws.Range("J138").Value = 149993.25
$ws.Range("L138").Value = 149993.25
$ws.Range("N138").Value = -160273.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 164.32259
$ws.Range("I7").Value = 76.86364
$ws.Range("J7").Value = 378.1111
$ws.Range("K7").Value = 76.86364
$ws.Range("L7").Value = 378.1111
$ws.Range("M7").Value = 36.13636
$ws.Range("N7").Value = -604.1111000000001

$ws.Range("H22").Value = 11769.667
$ws.Range("I22").Value = 25358
$ws.Range("J22").Value = 899
$ws.Range("K22").Value = 25358
$ws.Range("L22").Value = 899
$ws.Range("M22").Value = -25008
$ws.Range("N22").Value = -1599

$ws.Range("H31").Value = 3663.2173
$ws.Range("I31").Value = 2805.5454
$ws.Range("J31").Value = 4449.4165
$ws.Range("K31").Value = 2805.5454
$ws.Range("L31").Value = 4449.4165
$ws.Range("M31").Value = -2510.5454
$ws.Range("N31").Value = -5039.4165

$ws.Range("H34").Value = 3663.2173
$ws.Range("I34").Value = 2805.5454
$ws.Range("J34").Value = 4449.4165
$ws.Range("K34").Value = 2805.5454
$ws.Range("L34").Value = 4449.4165
$ws.Range("M34").Value = -2603.5454
$ws.Range("N34").Value = -4853.4165

$ws.Range("H58").Value = 5138.4287
$ws.Range("I58").Value = 5194.8335
$ws.Range("J58").Value = 4800
$ws.Range("K58").Value = 5194.8335
$ws.Range("L58").Value = 4800
$ws.Range("M58").Value = -4991.8335
$ws.Range("N58").Value = -5206

$ws.Range("H133").Value = 48130.2
$ws.Range("J133").Value = 60217
$ws.Range("L133").Value = 60217
$ws.Range("N133").Value = -65277

$ws.Range("H136").Value = 5138.4287
$ws.Range("I136").Value = 5194.8335
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 15584.5005
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -13034.5005
$ws.Range("N136").Value = -19500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2692.9167
$ws.Range("I131").Value = 911
$ws.Range("J131").Value = 3762.0667
$ws.Range("K131").Value = 2733
$ws.Range("L131").Value = 11286.2001
$ws.Range("M131").Value = 2307
$ws.Range("N131").Value = -21366.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 259.07144
$ws.Range("I2").Value = 111.72727
$ws.Range("J2").Value = 311.35483
$ws.Range("K2").Value = 111.72727
$ws.Range("L2").Value = 311.35483
$ws.Range("M2").Value = 1.272729999999996
$ws.Range("N2").Value = -537.35483

$ws.Range("H102").Value = 1437.1818
$ws.Range("I102").Value = 813.0294
$ws.Range("J102").Value = 3559.3
$ws.Range("K102").Value = 813.0294
$ws.Range("L102").Value = 3559.3
$ws.Range("M102").Value = 808.9706
$ws.Range("N102").Value = -6803.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1999.5
$ws.Range("I46").Value = 975.0968
$ws.Range("J46").Value = 4442.3076
$ws.Range("K46").Value = 975.0968
$ws.Range("L46").Value = 4442.3076
$ws.Range("M46").Value = -787.0968
$ws.Range("N46").Value = -4818.3076

$ws.Range("H82").Value = 2848.25
$ws.Range("I82").Value = 2225
$ws.Range("J82").Value = 4094.75
$ws.Range("K82").Value = 2225
$ws.Range("L82").Value = 4094.75
$ws.Range("M82").Value = -1864
$ws.Range("N82").Value = -4816.75

$ws.Range("H85").Value = 2848.25
$ws.Range("I85").Value = 2225
$ws.Range("J85").Value = 4094.75
$ws.Range("K85").Value = 2225
$ws.Range("L85").Value = 4094.75
$ws.Range("M85").Value = -977
$ws.Range("N85").Value = -6590.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3662.4546
$ws.Range("I132").Value = 3337.75
$ws.Range("J132").Value = 4528.3335
$ws.Range("K132").Value = 10013.25
$ws.Range("L132").Value = 13585.0005
$ws.Range("M132").Value = -7483.25
$ws.Range("N132").Value = -18645.0005

$ws.Range("H136").Value = 2220.7144
$ws.Range("I136").Value = 1442.4
$ws.Range("J136").Value = 4166.5
$ws.Range("K136").Value = 4327.200000000001
$ws.Range("L136").Value = 12499.5
$ws.Range("M136").Value = -1777.200000000001
$ws.Range("N136").Value = -17599.5
